$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date update
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting all subsequent rows (Description, Purpose, Copyright, Immutable) up by one.
$ws1.Rows.Item(11).Delete()
